# Natmi following Dr Hou advice
# Refresh the Nlgn2-Nrxn2 LR-pairs sheet: rows 2-3 (ECs sending cluster)
# get corrected values, and six additional rows (4-9) are added for the
# FAPs, M2 and sCs sending clusters, each paired with ECs/sCs targets.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,20
$row2[0,0] = "ECs"
$row2[0,1] = "Nlgn2"
$row2[0,2] = "Nrxn2"
$row2[0,3] = "ECs"
$row2[0,4] = 2
$row2[0,5] = 0.6666666666666666
$row2[0,6] = 2.164609666666667
$row2[0,7] = 6.493829
$row2[0,8] = 0.06870390863591093
$row2[0,9] = 0.06870390863591092
$row2[0,10] = 1
$row2[0,11] = 0.3333333333333333
$row2[0,12] = 0.0761
$row2[0,13] = 0.2283
$row2[0,14] = 0.3046644487415076
$row2[0,15] = 0.3046644487415076
$row2[0,16] = 0.1647267956333333
$row2[0,17] = 1.4825411607
$row2[0,18] = 0.02093163845094671
$row2[0,19] = 0.0209316384509467
$ws.Range("A2:T2").Value2 = $row2

$row3 = New-Object 'object[,]' 1,20
$row3[0,0] = "ECs"
$row3[0,1] = "Nlgn2"
$row3[0,2] = "Nrxn2"
$row3[0,3] = "sCs"
$row3[0,4] = 2
$row3[0,5] = 0.6666666666666666
$row3[0,6] = 2.164609666666667
$row3[0,7] = 6.493829
$row3[0,8] = 0.06870390863591093
$row3[0,9] = 0.06870390863591092
$row3[0,10] = 3
$row3[0,11] = 1
$row3[0,12] = 0.173683
$row3[0,13] = 0.521049
$row3[0,14] = 0.6953355512584923
$row3[0,15] = 0.6953355512584923
$row3[0,16] = 0.3759559007356667
$row3[0,17] = 3.383603106621
$row3[0,18] = 0.04777227018496422
$row3[0,19] = 0.04777227018496421
$ws.Range("A3:T3").Value2 = $row3

$row4 = New-Object 'object[,]' 1,20
$row4[0,0] = "FAPs"
$row4[0,1] = "Nlgn2"
$row4[0,2] = "Nrxn2"
$row4[0,3] = "ECs"
$row4[0,4] = 3
$row4[0,5] = 1
$row4[0,6] = 20.420946
$row4[0,7] = 61.262838
$row4[0,8] = 0.648153258228483
$row4[0,9] = 0.648153258228483
$row4[0,10] = 1
$row4[0,11] = 0.3333333333333333
$row4[0,12] = 0.0761
$row4[0,13] = 0.2283
$row4[0,14] = 0.3046644487415076
$row4[0,15] = 0.3046644487415076
$row4[0,16] = 1.5540339906
$row4[0,17] = 13.9863059154
$row4[0,18] = 0.1974692551181928
$row4[0,19] = 0.1974692551181928
$ws.Range("A4:T4").Value2 = $row4

$row5 = New-Object 'object[,]' 1,20
$row5[0,0] = "FAPs"
$row5[0,1] = "Nlgn2"
$row5[0,2] = "Nrxn2"
$row5[0,3] = "sCs"
$row5[0,4] = 3
$row5[0,5] = 1
$row5[0,6] = 20.420946
$row5[0,7] = 61.262838
$row5[0,8] = 0.648153258228483
$row5[0,9] = 0.648153258228483
$row5[0,10] = 3
$row5[0,11] = 1
$row5[0,12] = 0.173683
$row5[0,13] = 0.521049
$row5[0,14] = 0.6953355512584923
$row5[0,15] = 0.6953355512584923
$row5[0,16] = 3.546771164118
$row5[0,17] = 31.920940477062
$row5[0,18] = 0.4506840031102902
$row5[0,19] = 0.4506840031102902
$ws.Range("A5:T5").Value2 = $row5

$row6 = New-Object 'object[,]' 1,20
$row6[0,0] = "M2"
$row6[0,1] = "Nlgn2"
$row6[0,2] = "Nrxn2"
$row6[0,3] = "ECs"
$row6[0,4] = 3
$row6[0,5] = 1
$row6[0,6] = 0.3071453333333333
$row6[0,7] = 0.921436
$row6[0,8] = 0.009748679054813303
$row6[0,9] = 0.009748679054813303
$row6[0,10] = 1
$row6[0,11] = 0.3333333333333333
$row6[0,12] = 0.0761
$row6[0,13] = 0.2283
$row6[0,14] = 0.3046644487415076
$row6[0,15] = 0.3046644487415076
$row6[0,16] = 0.02337375986666667
$row6[0,17] = 0.2103638388
$row6[0,18] = 0.002970075930192577
$row6[0,19] = 0.002970075930192577
$ws.Range("A6:T6").Value2 = $row6

$row7 = New-Object 'object[,]' 1,20
$row7[0,0] = "M2"
$row7[0,1] = "Nlgn2"
$row7[0,2] = "Nrxn2"
$row7[0,3] = "sCs"
$row7[0,4] = 3
$row7[0,5] = 1
$row7[0,6] = 0.3071453333333333
$row7[0,7] = 0.921436
$row7[0,8] = 0.009748679054813303
$row7[0,9] = 0.009748679054813303
$row7[0,10] = 3
$row7[0,11] = 1
$row7[0,12] = 0.173683
$row7[0,13] = 0.521049
$row7[0,14] = 0.6953355512584923
$row7[0,15] = 0.6953355512584923
$row7[0,16] = 0.05334592292933334
$row7[0,17] = 0.480113306364
$row7[0,18] = 0.006778603124620726
$row7[0,19] = 0.006778603124620726
$ws.Range("A7:T7").Value2 = $row7

$row8 = New-Object 'object[,]' 1,20
$row8[0,0] = "sCs"
$row8[0,1] = "Nlgn2"
$row8[0,2] = "Nrxn2"
$row8[0,3] = "ECs"
$row8[0,4] = 3
$row8[0,5] = 1
$row8[0,6] = 8.613652999999999
$row8[0,7] = 25.840959
$row8[0,8] = 0.2733941540807927
$row8[0,9] = 0.2733941540807927
$row8[0,10] = 1
$row8[0,11] = 0.3333333333333333
$row8[0,12] = 0.0761
$row8[0,13] = 0.2283
$row8[0,14] = 0.3046644487415076
$row8[0,15] = 0.3046644487415076
$row8[0,16] = 0.6554989932999999
$row8[0,17] = 5.8994909397
$row8[0,18] = 0.08329347924217551
$row8[0,19] = 0.08329347924217551
$ws.Range("A8:T8").Value2 = $row8

$row9 = New-Object 'object[,]' 1,20
$row9[0,0] = "sCs"
$row9[0,1] = "Nlgn2"
$row9[0,2] = "Nrxn2"
$row9[0,3] = "sCs"
$row9[0,4] = 3
$row9[0,5] = 1
$row9[0,6] = 8.613652999999999
$row9[0,7] = 25.840959
$row9[0,8] = 0.2733941540807927
$row9[0,9] = 0.2733941540807927
$row9[0,10] = 3
$row9[0,11] = 1
$row9[0,12] = 0.173683
$row9[0,13] = 0.521049
$row9[0,14] = 0.6953355512584923
$row9[0,15] = 0.6953355512584923
$row9[0,16] = 1.496045093999
$row9[0,17] = 13.464405845991
$row9[0,18] = 0.1901006748386172
$row9[0,19] = 0.1901006748386172
$ws.Range("A9:T9").Value2 = $row9

